$wb = $excel.ActiveWorkbook

# --- ALC sheet updates ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J53").Value = 383.6
$ws.Range("M53").Value = 496.33333
$ws.Range("K53").Value = 140.66667
$ws.Range("L53").Value = 383.6
$ws.Range("I53").Value = 140.66667
$ws.Range("N53").Value = -1657.6
$ws.Range("H53").Value = 227.42857
$ws.Range("M64").Value = -3191.625
$ws.Range("K64").Value = 3439.625
$ws.Range("I64").Value = 3439.625
$ws.Range("H64").Value = 5842.136
$ws.Range("J64").Value = 12248.833
$ws.Range("N64").Value = -12744.833
$ws.Range("L64").Value = 12248.833
$ws.Range("K67").Value = 3439.625
$ws.Range("L67").Value = 12248.833
$ws.Range("N67").Value = -13964.833
$ws.Range("H67").Value = 5842.136
$ws.Range("J67").Value = 12248.833
$ws.Range("M67").Value = -2581.625
$ws.Range("I67").Value = 3439.625
$ws.Range("H74").Value = 10874.417
$ws.Range("L74").Value = 17998.8
$ws.Range("N74").Value = -19870.8
$ws.Range("J74").Value = 17998.8
$ws.Range("K74").Value = 5785.5713
$ws.Range("M74").Value = -4849.5713
$ws.Range("I74").Value = 5785.5713
$ws.Range("L77").Value = 89994
$ws.Range("K77").Value = 28927.8565
$ws.Range("N77").Value = -99354
$ws.Range("I77").Value = 5785.5713
$ws.Range("M77").Value = -24247.8565
$ws.Range("H77").Value = 10874.417
$ws.Range("J77").Value = 17998.8
$ws.Range("N97").Value = -167766.5
$ws.Range("J97").Value = 55591.5
$ws.Range("H97").Value = 37336
$ws.Range("L97").Value = 166774.5
$ws.Range("L99").Value = 1200
$ws.Range("H99").Value = 430.83334
$ws.Range("N99").Value = -4196
$ws.Range("J99").Value = 400
$ws.Range("I135").Value = 779.625
$ws.Range("M135").Value = -4481.625
$ws.Range("H135").Value = 1063.4286
$ws.Range("K135").Value = 7016.625
$ws.Range("K137").Value = 5233.5
$ws.Range("I137").Value = 1744.5
$ws.Range("H137").Value = 1859.1875
$ws.Range("M137").Value = -2683.5

# --- ARM sheet updates ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M61").Value = -2114.0908
$ws.Range("I61").Value = 2326.0908
$ws.Range("H61").Value = 3077.6428
$ws.Range("K61").Value = 2326.0908
$ws.Range("N107").Value = -57680
$ws.Range("L107").Value = 50000
$ws.Range("H107").Value = 50000
$ws.Range("J107").Value = 50000
$ws.Range("K110").Value = 10295.272
$ws.Range("H110").Value = 9115.5
$ws.Range("I110").Value = 10295.272
$ws.Range("M110").Value = -8250.272000000001
$ws.Range("I132").Value = 3128.1614
$ws.Range("L132").Value = 13777.0005
$ws.Range("J132").Value = 4592.3335
$ws.Range("M132").Value = -6854.484199999999
$ws.Range("N132").Value = -18837.0005
$ws.Range("H132").Value = 3257.353
$ws.Range("K132").Value = 9384.484199999999
$ws.Range("J135").Value = 60000
$ws.Range("L135").Value = 60000
$ws.Range("H135").Value = 60000
$ws.Range("N135").Value = -70140
$ws.Range("K136").Value = 6978.2724
$ws.Range("H136").Value = 3077.6428
$ws.Range("I136").Value = 2326.0908
$ws.Range("M136").Value = -4428.2724

# --- BSM sheet updates ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2552.647
$ws.Range("J86").Value = 6500
$ws.Range("K86").Value = 1338.0769
$ws.Range("N86").Value = -8746
$ws.Range("M86").Value = -215.0769
$ws.Range("L86").Value = 6500
$ws.Range("I86").Value = 1338.0769
$ws.Range("N89").Value = -43732
$ws.Range("K89").Value = 6690.3845
$ws.Range("I89").Value = 1338.0769
$ws.Range("L89").Value = 32500
$ws.Range("H89").Value = 2552.647
$ws.Range("M89").Value = -1074.3845
$ws.Range("J89").Value = 6500
$ws.Range("I94").Value = 969.7692
$ws.Range("K94").Value = 969.7692
$ws.Range("M94").Value = -518.7692
$ws.Range("L94").Value = 3960
$ws.Range("H94").Value = 1800.3889
$ws.Range("N94").Value = -4862
$ws.Range("J94").Value = 3960
$ws.Range("M134").Value = -2892.8799
$ws.Range("K134").Value = 5427.8799
$ws.Range("L134").Value = 7498.5
$ws.Range("I134").Value = 1809.2933
$ws.Range("J134").Value = 2499.5
$ws.Range("N134").Value = -12568.5
$ws.Range("H134").Value = 1827.2208
$ws.Range("J135").Value = 89000
$ws.Range("L135").Value = 89000
$ws.Range("H135").Value = 89000
$ws.Range("N135").Value = -99140

# --- CRP sheet updates ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H64").Value = 12500
$ws.Range("J64").Value = 12500
$ws.Range("N64").Value = -12996
$ws.Range("L64").Value = 12500
$ws.Range("L67").Value = 12500
$ws.Range("N67").Value = -14216
$ws.Range("H67").Value = 12500
$ws.Range("J67").Value = 12500
$ws.Range("L105").Value = 1437.5
$ws.Range("J105").Value = 1437.5
$ws.Range("K105").Value = 1274.9166
$ws.Range("I105").Value = 1274.9166
$ws.Range("H105").Value = 1315.5625
$ws.Range("M105").Value = 472.0834
$ws.Range("N105").Value = -4931.5
$ws.Range("I132").Value = 4541.222
$ws.Range("L132").Value = 11999.4
$ws.Range("J132").Value = 3999.8
$ws.Range("M132").Value = -11093.666
$ws.Range("N132").Value = -17059.4
$ws.Range("H132").Value = 4347.857
$ws.Range("K132").Value = 13623.666
$ws.Range("M134").Value = -19543.7139
$ws.Range("K134").Value = 22078.7139
$ws.Range("I134").Value = 7359.5713
$ws.Range("H134").Value = 13388.682
$ws.Range("K138").Value = 94999
$ws.Range("I138").Value = 94999
$ws.Range("H138").Value = 94999
$ws.Range("M138").Value = -89859

# --- CUL sheet updates ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M26").Value = -1294.125
$ws.Range("H26").Value = 1130.4445
$ws.Range("I26").Value = 527.375
$ws.Range("K26").Value = 1582.125

# --- GSM sheet updates ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M19").Value = -19712
$ws.Range("K19").Value = 20000
$ws.Range("H19").Value = 13247.5
$ws.Range("I19").Value = 20000
$ws.Range("L19").Value = 10996.667
$ws.Range("J19").Value = 10996.667
$ws.Range("N19").Value = -11572.667
$ws.Range("I70").Value = 5788.125
$ws.Range("K70").Value = 5788.125
$ws.Range("M70").Value = -5518.125
$ws.Range("N70").Value = -26541
$ws.Range("J70").Value = 26001
$ws.Range("H70").Value = 8034
$ws.Range("L70").Value = 26001
$ws.Range("M73").Value = -4852.125
$ws.Range("K73").Value = 5788.125
$ws.Range("I73").Value = 5788.125
$ws.Range("J73").Value = 26001
$ws.Range("L73").Value = 26001
$ws.Range("N73").Value = -27873
$ws.Range("H73").Value = 8034
$ws.Range("I132").Value = 3323.375
$ws.Range("M132").Value = -7440.125
$ws.Range("K132").Value = 9970.125
$ws.Range("H132").Value = 4197.0835

# --- LTW sheet updates ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("L20").Value = 15000
$ws.Range("H20").Value = 15000
$ws.Range("N20").Value = -15452
$ws.Range("J20").Value = 15000
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("L64").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("M68").Value = -1480.6667
$ws.Range("K68").Value = 2229.6667
$ws.Range("I68").Value = 2229.6667
$ws.Range("H68").Value = 324661.25
$ws.Range("H71").Value = 324661.25
$ws.Range("I71").Value = 2229.6667
$ws.Range("M71").Value = -7404.333500000001
$ws.Range("K71").Value = 11148.3335
$ws.Range("J75").Value = 37500
$ws.Range("H75").Value = 37500
$ws.Range("L75").Value = 37500
$ws.Range("N75").Value = -39372
$ws.Range("H78").Value = 37500
$ws.Range("L78").Value = 112500
$ws.Range("J78").Value = 37500
$ws.Range("N78").Value = -121860
$ws.Range("M100").Value = -10583
$ws.Range("J100").Value = 26249.5
$ws.Range("I100").Value = 11124
$ws.Range("K100").Value = 11124
$ws.Range("N100").Value = -27331.5
$ws.Range("L100").Value = 26249.5
$ws.Range("H100").Value = 16165.833
$ws.Range("I132").Value = 2386.52
$ws.Range("M132").Value = -4629.559999999999
$ws.Range("K132").Value = 7159.559999999999
$ws.Range("H132").Value = 3077.3547
$ws.Range("L133").Value = 66835.14
$ws.Range("N133").Value = -71895.14
$ws.Range("H133").Value = 66835.14
$ws.Range("J133").Value = 66835.14

# --- WVR sheet updates ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 1000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("N70").Value = -40735
$ws.Range("J70").Value = 40105
$ws.Range("H70").Value = 40105
$ws.Range("L70").Value = 40105
$ws.Range("J73").Value = 40105
$ws.Range("L73").Value = 40105
$ws.Range("N73").Value = -42289
$ws.Range("H73").Value = 40105
$ws.Range("M100").Value = -869
$ws.Range("J100").Value = 708.75
$ws.Range("I100").Value = 705
$ws.Range("K100").Value = 1410
$ws.Range("N100").Value = -2499.5
$ws.Range("L100").Value = 1417.5
$ws.Range("H100").Value = 706.0345
$ws.Range("K136").Value = 4051.1613
$ws.Range("H136").Value = 1517.7894
$ws.Range("J136").Value = 2259.1428
$ws.Range("I136").Value = 1350.3871
$ws.Range("M136").Value = -1501.1613
$ws.Range("L136").Value = 6777.428400000001
$ws.Range("N136").Value = -11877.4284
